$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1821.0769
$ws.Range("I125").Value = 474.5
$ws.Range("J125").Value = 2419.5557
$ws.Range("K125").Value = 4270.5
$ws.Range("L125").Value = 21776.0013
$ws.Range("M125").Value = -1810.5
$ws.Range("N125").Value = -26696.0013

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 35049
$ws.Range("J44").Value = 35049
$ws.Range("L44").Value = 35049
$ws.Range("N44").Value = -36025

$ws.Range("H55").Value = 21473.334
$ws.Range("J55").Value = 21473.334
$ws.Range("L55").Value = 21473.334
$ws.Range("N55").Value = -22103.334

$ws.Range("H80").Value = 24860.6
$ws.Range("J80").Value = 24860.6
$ws.Range("L80").Value = 24860.6
$ws.Range("N80").Value = -26856.6

$ws.Range("H83").Value = 24860.6
$ws.Range("J83").Value = 24860.6
$ws.Range("L83").Value = 74581.79999999999
$ws.Range("N83").Value = -84565.79999999999

$ws.Range("H122").Value = 2350.2222
$ws.Range("J122").Value = 2700
$ws.Range("L122").Value = 8100
$ws.Range("N122").Value = -13000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34987
$ws.Range("J35").Value = 34987
$ws.Range("L35").Value = 34987
$ws.Range("N35").Value = -35607

$ws.Range("H82").Value = 48015.04
$ws.Range("J82").Value = 29430.79
$ws.Range("L82").Value = 29430.79
$ws.Range("N82").Value = -30196.79

$ws.Range("H85").Value = 48015.04
$ws.Range("J85").Value = 29430.79
$ws.Range("L85").Value = 29430.79
$ws.Range("N85").Value = -32082.79

$ws.Range("H122").Value = 29621.818
$ws.Range("J122").Value = 29621.818
$ws.Range("L122").Value = 29621.818
$ws.Range("N122").Value = -39421.818

$ws.Range("H132").Value = 49960
$ws.Range("J132").Value = 49960
$ws.Range("L132").Value = 49960
$ws.Range("N132").Value = -60080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 303956
$ws.Range("J97").Value = 303956
$ws.Range("L97").Value = 303956
$ws.Range("N97").Value = -305938

$ws.Range("H99").Value = 1784.8889
$ws.Range("I99").Value = 1635.6
$ws.Range("J99").Value = 1971.5
$ws.Range("K99").Value = 1635.6
$ws.Range("L99").Value = 1971.5
$ws.Range("M99").Value = -137.5999999999999
$ws.Range("N99").Value = -4967.5

$ws.Range("H107").Value = 1076.3684
$ws.Range("I107").Value = 1086.0588
$ws.Range("J107").Value = 994
$ws.Range("K107").Value = 1086.0588
$ws.Range("L107").Value = 994
$ws.Range("M107").Value = 833.9412
$ws.Range("N107").Value = -4834

$ws.Range("H126").Value = 1784.8889
$ws.Range("I126").Value = 1635.6
$ws.Range("J126").Value = 1971.5
$ws.Range("K126").Value = 4906.799999999999
$ws.Range("L126").Value = 5914.5
$ws.Range("M126").Value = -2436.799999999999
$ws.Range("N126").Value = -10854.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 3403.3333
$ws.Range("I56").Value = 3403.3333
$ws.Range("K56").Value = 3403.3333
$ws.Range("M56").Value = -2873.3333

$ws.Range("H112").Value = 36277230
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 43914336
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 131743008
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -131745224

$ws.Range("H122").Value = 3531.3562
$ws.Range("I122").Value = 407.14285
$ws.Range("J122").Value = 3804.725
$ws.Range("K122").Value = 3664.28565
$ws.Range("L122").Value = 34242.525
$ws.Range("M122").Value = -1214.28565
$ws.Range("N122").Value = -39142.525

$ws.Range("H125").Value = 2431.818
$ws.Range("J125").Value = 2476.1904
$ws.Range("L125").Value = 7428.5712
$ws.Range("N125").Value = -17268.5712

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2306.0476
$ws.Range("I102").Value = 1653.6154
$ws.Range("J102").Value = 3366.25
$ws.Range("K102").Value = 1653.6154
$ws.Range("L102").Value = 3366.25
$ws.Range("M102").Value = -31.61539999999991
$ws.Range("N102").Value = -6610.25

$ws.Range("H122").Value = 2357.1428
$ws.Range("I122").Value = 2750
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 8250
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -5800
$ws.Range("N122").Value = -10400.0002

$ws.Range("H123").Value = 10203.75
$ws.Range("J123").Value = 10203.75
$ws.Range("L123").Value = 10203.75
$ws.Range("N123").Value = -15103.75

$ws.Range("H126").Value = 4869.2188
$ws.Range("I126").Value = 2503.6428
$ws.Range("J126").Value = 6709.1113
$ws.Range("K126").Value = 7510.928400000001
$ws.Range("L126").Value = 20127.3339
$ws.Range("M126").Value = -5040.928400000001
$ws.Range("N126").Value = -25067.3339

$ws.Range("H128").Value = 50447.145
$ws.Range("J128").Value = 50447.145
$ws.Range("L128").Value = 50447.145
$ws.Range("N128").Value = -60407.145

$ws.Range("H132").Value = 3185.7407
$ws.Range("I132").Value = 2640.9333
$ws.Range("J132").Value = 3866.75
$ws.Range("K132").Value = 7922.7999
$ws.Range("L132").Value = 11600.25
$ws.Range("M132").Value = -5392.7999
$ws.Range("N132").Value = -16660.25

$ws.Range("H136").Value = 32040.75
$ws.Range("J136").Value = 32040.75
$ws.Range("L136").Value = 96122.25
$ws.Range("N136").Value = -101222.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6668314.5
$ws.Range("I7").Value = 7144336.5
$ws.Range("K7").Value = 7144336.5
$ws.Range("M7").Value = -7144224.5

$ws.Range("H40").Value = 1470.5238
$ws.Range("I40").Value = 1419.05
$ws.Range("K40").Value = 1419.05
$ws.Range("M40").Value = -1283.05

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 2377.7856
$ws.Range("I122").Value = 1924.9166
$ws.Range("K122").Value = 5774.7498
$ws.Range("M122").Value = -3324.7498

$ws.Range("H126").Value = 6668314.5
$ws.Range("I126").Value = 7144336.5
$ws.Range("K126").Value = 21433009.5
$ws.Range("M126").Value = -21430539.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12361.6
$ws.Range("I122").Value = 18751
$ws.Range("K122").Value = 56253
$ws.Range("M122").Value = -53803

$ws.Range("H126").Value = 2671.6428
$ws.Range("I126").Value = 1675.375
$ws.Range("K126").Value = 5026.125
$ws.Range("M126").Value = -2556.125
